$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "dsadas"
$ws.Range("B2").Value = 0

$ws.Range("A3").Value = $null
$ws.Range("B3").Value = $null
